$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("WC49 P5H", "Screw K30 no lo detecta puesto", "2024-06-03", "12:23:04", "Mañana", "12:23:06", "0:00:02", "N/A"),
    @("WC49 P5H", "La cámara no detecta Top Cover", "2024-06-03", "12:23:11", "Mañana", "12:23:13", "0:00:02", "N/A"),
    @("WC49 P5H", "Screw K30 no lo detecta puesto", "2024-06-03", "12:23:16", "Mañana", "12:23:17", "0:00:01", "0.12 minutos"),
    @("WC48 P5F", "AOI (fallo etiqueta)",            "2024-06-03", "12:39:40", "Mañana", "12:39:42", "0:00:02", "N/A"),
    @("WC48 P5F", "Etiquetadora",                    "2024-06-03", "12:39:49", "Mañana", "12:39:50", "0:00:01", "N/A"),
    @("WC48 P5F", "Etiquetadora",                    "2024-06-03", "12:39:51", "Mañana", "12:39:52", "0:00:01", "0.14 minutos"),
    @("WC48 P5F", "Etiquetadora",                    "2024-06-03", "12:39:53", "Mañana", "12:39:54", "0:00:01", "0.09 minutos")
)

$startRow = 110
$endRow = $startRow + $rows.Count - 1
$ws.Range("C$startRow`:C$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
